# Rows 28 and 29 of the active sheet represent two species observations
# that were swapped (re-ordered) in the source data. Update each changed
# cell in place so that row 28 now holds what used to be row 29's data
# (plus the empty "Kön" cell in column L) and row 29 now holds what used
# to be row 28's data (and loses the column L cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 28: becomes the former row 29 content ----
$ws.Range("A28").Value = 111596897
$ws.Range("B28").Value = 103288
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 221144
$ws.Range("F28").Value = "Grönpyrola"
$ws.Range("G28").Value = "Pyrola chlorantha"
$ws.Range("H28").Value = "Sw."
$ws.Range("I28").Value = "'100"
$ws.Range("J28").Value = "plantor/tuvor"
# Materialize an (empty) L28 cell, mirroring the neighboring K28/N28
# placeholder cells, since the diff adds a blank <c r="L28"/> cell here.
$ws.Range("L28").Font.Bold = $ws.Range("L28").Font.Bold
$ws.Range("P28").Value = "Björkmossen 227 m E, Upl"
$ws.Range("Q28").Value = 654422.181084068
$ws.Range("R28").Value = 6690769.97221576
$ws.Range("Z28").Value = "12:53"
$ws.Range("AB28").Value = "12:53"
$ws.Range("AC28").Value = "Uppskattat antal, helt tjockt med plantor så går ej att räkna."

# ---- Row 29: becomes the former row 28 content ----
$ws.Range("A29").Value = 111596843
$ws.Range("B29").Value = 90709
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 5448
$ws.Range("F29").Value = "Svartvit taggsvamp"
$ws.Range("G29").Value = "Phellodon connatus"
$ws.Range("H29").Value = "(Schultz) nom.prov"
$ws.Range("I29").Value = "'20"
$ws.Range("J29").Value = "fruktkroppar"
$ws.Range("L29").ClearContents()
$ws.Range("P29").Value = "Björkmossen 282 m E, Upl"
$ws.Range("Q29").Value = 654476.3214109741
$ws.Range("R29").Value = 6690758.38440035
$ws.Range("Z29").Value = "12:00"
$ws.Range("AB29").Value = "12:02"
$ws.Range("AC29").Value = "Väldigt svag doft."
